$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "県"
$ws.Range("A2").Value = "福岡"
$ws.Range("A3").Value = "佐賀"
$ws.Range("A4").Value = "長崎"
$ws.Range("A5").Value = "熊本"
$ws.Range("A6").Value = "大分"
$ws.Range("A7").Value = "宮崎"
$ws.Range("A8").Value = "鹿児島"

$ws.Range("A8").Select()
